$d = $word.ActiveDocument

# Hybrid bold + color highlight used for quantitative impact metrics
# (percentages, dollar amounts, large numbers).
# Target OOXML color is RGB 2C3E50; the Word object model Font.Color
# value is packed as 0x00BBGGRR, i.e. B=0x50,G=0x3E,R=0x2C -> 0x503E2C.
$metricColor = 5258796

function Get-ParagraphByText($uniqueSubstring, $excludeSubstring) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $t = $p.Range.Text
        if ($t -like "*$uniqueSubstring*") {
            if ($excludeSubstring -ne $null -and $t -like "*$excludeSubstring*") {
                continue
            }
            return $p
        }
    }
    Write-Output "PARAGRAPH NOT FOUND: $uniqueSubstring"
    return $null
}

function Highlight-Metric($paragraph, $metricText) {
    if ($paragraph -eq $null) { return }
    $rng = $paragraph.Range
    $found = $rng.Find.Execute($metricText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "METRIC NOT FOUND: $metricText"
        return
    }
    $rng.Font.Bold = 1
    $rng.Font.Color = $metricColor
}

# 1) "Discovered systematic race coding errors ... from 23% to 64%"
$para1 = Get-ParagraphByText "developed geospatial machine learning algorithms improving demographic classification accuracy" $null
Highlight-Metric $para1 "23%"
Highlight-Metric $para1 "64%"

# 2) "Achieved 87% prediction accuracy ... ±4.2% to ±2.1%" (long bullet, under Siege Analytics)
$plusMinus = [char]0x00B1
$metric42 = $plusMinus + "4.2%"
$metric21 = $plusMinus + "2.1%"
$para2 = Get-ParagraphByText "reducing polling error margins" $null
Highlight-Metric $para2 "87%"
Highlight-Metric $para2 "71%"
Highlight-Metric $para2 $metric42
Highlight-Metric $para2 $metric21

# 3) "Wrote RFP and analyzed bids from 1,200 vendors ..."
$para3 = Get-ParagraphByText "Wrote RFP and analyzed bids from" $null
Highlight-Metric $para3 "1,200"

# 4) "Created comprehensive meta-analysis framework ... $400M ... $1B+"
$para4 = Get-ParagraphByText "Created comprehensive meta-analysis framework" $null
Highlight-Metric $para4 '$400M'
Highlight-Metric $para4 '$1B'

# 5) "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
$para5 = Get-ParagraphByText "Algorithm reduced mapping costs by" $null
Highlight-Metric $para5 "73.5%"
Highlight-Metric $para5 '$4.7M'

# 6) "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
#    (short bullet, Key Achievements section — exclude the longer bullet from #2)
$para6 = Get-ParagraphByText "Achieved 87%" "reducing polling error margins"
Highlight-Metric $para6 "87%"
Highlight-Metric $para6 "71%"

Write-Output "Done"
